# "ajout save to excel" - populate the previously empty FDK totals
# (columns I:M, rows 4-9) on the "Scores" sheet with the computed
# Glenoid totals, and restore the last active selection (L16).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scores")

# Values for I4:M9 (5 columns x 6 rows), taken from the target workbook.
$data = @(
    @(3852.6358305144399, 3886.8614198381401, 5711.1949775252497, 7132.6198831584798, 7670.3851424714703),
    @(4302.4431323520403, 3168.0107132970802, 5613.5523332412704, 7733.9349808274601, 8308.1378704456602),
    @(5149.4458870725002, 3226.2202383661602, 5323.1651638346402, 8262.5852476171094, 8889.24394139182),
    @(6223.6801328150696, 3551.4517577358401, 4869.8091767357,    8840.2381634644407, 9564.2249611155694),
    @(8216.3184893125908, 4896.2539788139602, 5180.4819909731204, 9922.4508294343195, 10140.121829494499),
    @(8712.2991309537902, 6101.0087343801097, 6567.5330099508801, 10983.1622659615,   10980.1522210362)
)

$rowCount = $data.Length
$colCount = $data[0].Length

# Excel expects a 2-dimensional array (rows x columns) when assigning
# .Value to a multi-cell range.
$arr = New-Object 'object[,]' $rowCount,$colCount
for ($r = 0; $r -lt $rowCount; $r++) {
    for ($c = 0; $c -lt $colCount; $c++) {
        $arr[$r, $c] = $data[$r][$c]
    }
}

$ws.Range("I4:M9").Value = $arr

# Restore the cell selection left by the author (was G5, now L16).
$ws.Range("L16").Select()
